# Remove the trailing "Cohort" column from the CasesTab Neo4j query (B2),
# matching the upstream commit that trimmed the query text, and refresh the
# row height / selection state that Excel re-records when a cell with
# wrapped text shrinks by one line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the CasesTab query text in B2 (drop the trailing cohort clause).
$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n`nMATCH (c)<--(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

# 2. The query text is now one wrapped line shorter (21 -> 20 lines at
#    15pt/line), so the row needs to shrink from 315 to 300 points.
$ws.Rows.Item(2).RowHeight = 300

# 3. Move the active selection to B2 (matches the saved view state).
$ws.Range("B2").Select()

# 4. Best-effort: restore the scrolled viewport to row 2 as well.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
